# Revert "Use different months for summer for bsth, os, cs"
# Collapse the per-application-method Summer rows on the Climate sheet back
# down to a single Summer row, and drop the now-unused app.mthd.wthr /
# notes columns. Also restores the previously-active sheet/selection.

$wb = $excel.ActiveWorkbook

$wsClimate = $wb.Worksheets.Item("Climate")

# Drop the two extra "Summer" rows that had been split out per application
# method (open slot injection / closed slot injection). Row 5 (trailing
# hose / generic summer) is kept and becomes the single Summer row.
$wsClimate.Rows.Item(8).Delete() | Out-Null
$wsClimate.Rows.Item(7).Delete() | Out-Null

# Drop the "notes" column (last column, G) and the "app.mthd.wthr" column
# (C) which only existed to differentiate the per-method summer rows.
$wsClimate.Columns.Item(7).Delete() | Out-Null
$wsClimate.Columns.Item(3).Delete() | Out-Null

# Restore the previous view state: Climate was selected at B9 before this
# change set; after reverting, the selection on Climate moves to E16 and
# Application becomes the active/selected tab at B7.
$wsClimate.Range("E16").Select() | Out-Null

$wsApplication = $wb.Worksheets.Item("Application")
$wsApplication.Range("B7").Select() | Out-Null
$wsApplication.Activate() | Out-Null
